$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update prices for STM32F411RE (row 3) and wire kit (row 4)
$ws.Range("C3").Value = 391
$ws.Range("C4").Value = 304

# Mark both as obtained ("SI")
$ws.Range("D3").Value = "SI"
$ws.Range("D4").Value = "SI"

# Copy the alignment/number-format style from D9 (an existing "SI" cell) onto D3/D4
$ws.Range("D9").Copy()
$ws.Range("D3:D4").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Add a new formatted (underlined) empty cell at H8 and select it,
# matching the stray formatted cell that appears in the edited sheet
$ws.Range("H8").Style = "Normal"
$ws.Range("H8").Font.Underline = $true
$ws.Range("H8").Select()
